$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 0) Insert two new rows before the footnote row (currently row 16)
#    to host the new "ingreso" balance row (mean/SE line + SE-only
#    line). This pushes the footnote row down from row 16 to 18.
# ---------------------------------------------------------------
$ws.Rows("16:17").Insert()

# Work area that will receive text (not auto-numeric) values.
$textRange = $ws.Range("A1:G18")
$textRange.NumberFormat = "@"

# ---------------------------------------------------------------
# 1) Update existing Mean/(SE) and Mean-difference figures that
#    changed because the iebaltab summary was regenerated with the
#    new "ingreso" variable included in the balance table.
# ---------------------------------------------------------------
$ws.Range("C4").Value = "0.068"
$ws.Range("E4").Value = "0.082"
$ws.Range("G4").Value = "0.014"
$ws.Range("C5").Value = "(0.009)"
$ws.Range("E5").Value = "(0.035)"

$ws.Range("C6").Value = "0.158"
$ws.Range("E6").Value = "0.393"
$ws.Range("G6").Value = "0.235***"

$ws.Range("C8").Value = "0.074"
$ws.Range("G8").Value = "0.073**"

$ws.Range("C10").Value = "0.144"
$ws.Range("E10").Value = "0.197"
$ws.Range("G10").Value = "0.053"
$ws.Range("C11").Value = "(0.031)"
$ws.Range("E11").Value = "(0.112)"

$ws.Range("C12").Value = "0.202"
$ws.Range("E12").Value = "0.508"
$ws.Range("G12").Value = "0.306***"
$ws.Range("C13").Value = "(0.019)"
$ws.Range("E13").Value = "(0.089)"

$ws.Range("C14").Value = "0.089"
$ws.Range("G14").Value = "0.108**"
$ws.Range("C15").Value = "(0.012)"

# ---------------------------------------------------------------
# 2) Fill the new "ingreso" balance row (835 vs 61 obs, means,
#    standard errors and the pairwise mean-difference figure).
# ---------------------------------------------------------------
$ws.Range("A16").Value = "ingreso"
$ws.Range("B16").Value = "835"
$ws.Range("C16").Value = "2.792"
$ws.Range("D16").Value = "61"
$ws.Range("E16").Value = "3.660"
$ws.Range("F16").Value = "896"
$ws.Range("G16").Value = "0.868***"

$ws.Range("A17").Value = ""
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = "(0.052)"
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = "(0.092)"
$ws.Range("F17").Value = ""
$ws.Range("G17").Value = ""

# ---------------------------------------------------------------
# 3) Update the footnote text (now on row 18) so the "Full user
#    input" note mentions the new "ingreso" variable too.
# ---------------------------------------------------------------
$ws.Range("A18").Value = "If the table includes missing values (.n, .o, .v etc.) see the Missing values section in the help file for the Stata command iebaltab for definitions of these values. Significance: ***=.01, **=.05, *=.1. Full user input as written by user: [iebaltab dummy_jb dummy_d1 dummy_ara cantidad_jb cantidad_d1 cantidad_ara ingreso , groupvar(dummy_oxxo) control(0) savexlsx(difmedias_controles_staggered_variables_2019) replace] "

# Drop the temporary text formatting again so the saved styles stay
# as close as possible to the original (values remain text/strings
# because Excel keeps the cached cell type once entered as text).
$ws.Range("A1:G18").ClearFormats()
